# household_new.xlsx — "changed to a different directory"
#
# Content-level changes (once the shared-string garbage collection that
# happens on save is accounted for, every other numeric shift in the
# original diff is just a re-index of otherwise-unchanged strings):
#
#   1. The "queries" sheet's auxillaryHash formula-string literal changes
#      from escape(...) to encodeURIComponent(...) in G2:G4.
#   2. G4's stray bold/quote-prefix style is normalized to match G2:G3.
#   3. The active tab moves from "household" to "queries" (and the
#      household sheet's selection/tab state reverts to not-selected).

$wb = $excel.ActiveWorkbook

$queries = $wb.Worksheets.Item("queries")
$household = $wb.Worksheets.Item("household")

# A leading single quote in a .Value assignment is interpreted as Excel's
# "treat the rest as text" quote-prefix marker and gets swallowed, so the
# literal leading apostrophe in this formula-like string literal has to be
# doubled up to survive the round trip.
$newFormula = "''household_id='+encodeURIComponent(data('household_id'))"
$queries.Range("G2").Value = $newFormula
$queries.Range("G3").Value = $newFormula
$queries.Range("G4").Value = $newFormula

# G4 previously carried a one-off bold + quotePrefix style; line it up with
# G2/G3's style so it stops being a unique outlier.
$queries.Range("G4").NumberFormat = $queries.Range("G2").NumberFormat

# household was the active/selected tab; queries becomes active instead,
# with its selection on G5 (was D24). household keeps its own selection
# (H7) but loses tabSelected.
$household.Activate()
$queries.Activate()
$queries.Range("G5").Select()
